$d = $word.ActiveDocument

# --- Step 1: add strike+highlight formatting to the four plain paragraphs ---
# (paragraph indices are stable across these replacements since each InsertXML
#  call replaces exactly one paragraph's content with exactly one new paragraph)

$p2Xml = @'
<w:p w14:paraId="73DC8482" w14:textId="2B51496F" w:rsidR="00842681" w:rsidRPr="004E7D55" w:rsidRDefault="00034E47" w:rsidP="00034E47"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="004E7D55"><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>Įdėt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="004E7D55"><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="004E7D55"><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>Zenerio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="004E7D55"><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> diodus, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="004E7D55"><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>viršįtampiams</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="004E7D55"><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> ant L9637 VS ir VCC pinų;</w:t></w:r></w:p>
'@
$d.Paragraphs(2).Range.InsertXML($p2Xml)

$p5Xml = @'
<w:p w14:paraId="4A07EEF5" w14:textId="15BB63DD" w:rsidR="00EE32FE" w:rsidRDefault="00EE32FE" w:rsidP="00034E47"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>Įdėt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> USB </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>micro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> arba USB-C, kuris bus naudojamas </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>užmaitint</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>boardą</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">, kai </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>bootloadinsim</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>;</w:t></w:r></w:p>
'@
$d.Paragraphs(5).Range.InsertXML($p5Xml)

$p9Xml = @'
<w:p w14:paraId="62B1906E" w14:textId="6D7818E1" w:rsidR="00521AC9" w:rsidRDefault="00521AC9" w:rsidP="00034E47"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>Pakeisti L7805CV konverteri į SMD tipo;</w:t></w:r></w:p>
'@
$d.Paragraphs(9).Range.InsertXML($p9Xml)

$p12Xml = @'
<w:p w14:paraId="43903C0D" w14:textId="05C0C8D3" w:rsidR="00521AC9" w:rsidRDefault="00521AC9" w:rsidP="00034E47"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Visus komponentus </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>naudot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> SMD ir </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>sukelt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> į bendras projekto bibliotekas;</w:t></w:r></w:p>
'@
$d.Paragraphs(12).Range.InsertXML($p12Xml)

# --- Step 2: append the new paragraphs (Capacitoriai / parametrai / component list) ---
# Insert at a collapsed range right after the last paragraph's end, so nothing existing
# gets overwritten and the new paragraphs land right before the sectPr, as in the diff.

$newParasXml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>Capacitoriai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> surasti:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">100 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>nF</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">10 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>uF</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">330 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>nF</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>Surasyti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> parametrus:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>RES 100</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>RES 1K</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>RES 510</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>CAP 100nF</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>CAP 330nF</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>CAP 1uF</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>CAP 10uF</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>LED</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>ZENON 5V</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>ZENON 20V</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>STM32G0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>Connector</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>L7805ABD2T-TR</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>LM1117MP-3.3/NOPB</w:t></w:r><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>L9637D013TR</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>BM20</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>DIODAI</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>MYGTUKAS</w:t></w:r></w:p>
'@

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$insertPoint.InsertXML($newParasXml)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
